$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at C and D (old C/D shift right to E/F)
$ws.Columns.Item(3).Insert()
$ws.Columns.Item(4).Insert()

# Header row (bold, centered, no border) for the two new columns
$ws.Range("C1").Value = "Doc"
$ws.Range("D1").Value = "Type"
$ws.Range("C1:D1").Font.Bold = $true
$ws.Range("C1:D1").HorizontalAlignment = -4108

# Fill Doc/Type values for the three document-type blocks
$ws.Range("C2:C92").Value = "EPANB"
$ws.Range("D2:D92").Value = "MNB"

$ws.Range("C93:C136").Value = "CDN"
$ws.Range("D93:D136").Value = "Metas de las CDN"

$ws.Range("C137:C153").Value = "E50"
$ws.Range("D137:D153").Value = "Otras metas"

# Column widths
$ws.Columns.Item(2).ColumnWidth = 25.33203125
$ws.Columns.Item(3).ColumnWidth = 6.1640625
$ws.Columns.Item(4).ColumnWidth = 13.83203125
$ws.Columns.Item(5).ColumnWidth = 9.83203125
$ws.Columns.Item(6).ColumnWidth = 21.1640625
